$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: fix the duplicated "the  The" typo into
# "the supporting blocks.) The" so the sentence reads cleanly.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("of the  The springs", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $rng1.Find.Found) {
    throw "Could not find target text for edit 1"
}
$rng1.Text = "of the supporting blocks.) The springs"

# Force a run split so "supporting blocks.)" becomes its own run,
# matching how the author re-typed just that phrase.
$rngBlocks = $d.Content
$rngBlocks.Find.Execute("supporting blocks.)") | Out-Null
$rngBlocks.Font.Bold = 1
$rngBlocks.Font.Bold = 0

# ---------------------------------------------------------------------
# Edit 2: rewrite the spring-force justification and add two new
# paragraphs of discussion about the second spring / zipline clearance.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("provide a reasonable angular acceleration") | Out-Null
if (-not $rng2.Find.Found) {
    throw "Could not find target text for edit 2a"
}
$rng2.Text = "allow the mechanism to spring back into place quickly"

# Force a run split so the new phrase is its own run (boundaries at
# both ends of the replaced range already coincide with the adjacent
# runs, so this creates the 3-way split seen in the diff).
$rngSpring = $d.Content
$rngSpring.Find.Execute("allow the mechanism to spring back into place quickly") | Out-Null
$rngSpring.Font.Bold = 1
$rngSpring.Font.Bold = 0

# Add the two new paragraphs after the (current) last paragraph.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$para2 = $d.Paragraphs.Last
$para2.Range.Text = "The force on the second spring doesn" + [char]0x2019 + "t need to be very large to counteract rotation when the basket is on the zipline, as the perpendicular distance from the spring to the lower rotation axis is much larger than the distance for the normal force applied by the rod. Since this normal force is on the order of 100 N, we can estimate that the minimum required force for the second spring is about 20 to 30 N, although higher forces will allow the mechanism to spring back faster."

$para2.Range.InsertParagraphAfter()

$para3 = $d.Paragraphs.Last
$para3.Range.Text = "It is also important that the springs not be too strong, or the lifting mechanism will not be able to push the supports through the zipline. This will be more easily tuned using real springs once the assembly is built, but the forces seem to be small enough to be manageable with a DC motor."

# Move the _GoBack bookmark to the end of the new final paragraph (right
# before its paragraph mark), where Word would have left it after the
# author's last keystroke. Adding a bookmark exactly at the current
# document end confuses the host engine, so temporarily pad the
# document by one character while we reposition it, then trim the pad.
if ($d.Bookmarks.Exists("_GoBack")) {
    $finalParaRange = $d.Paragraphs.Last.Range
    $endPos = $finalParaRange.End - 1

    $pad = $d.Range($endPos, $endPos)
    $pad.InsertAfter("X")

    $targetRange = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("_GoBack", $targetRange)

    $padRange = $d.Range($endPos, $endPos + 1)
    $padRange.Text = ""
}

Write-Output "Done"
